# Update column G ("K" - strikeouts) values on the active sheet to reflect
# the regenerated save_data (K instead of Strike#, regen std/mean, calc and
# write s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    8  = 0
    9  = 1
    10 = 0
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    17 = 1
    18 = 1
    19 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
